$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    paragraph at the top of the document.
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaTail = ": Discover the top features of 1 Left Alive slot game, including high payout percentage and special symbols. Play for free or real money."
$metaXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
    "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>$metaTail</w:t></w:r></w:p></w:body>" +
    "</w:document></pkg:xmlData></pkg:part></pkg:package>"

$metaPara.Range.InsertXML($metaXml)

# -----------------------------------------------------------------
# 2) At the very end of the document, drop the duplicated bold title
#    paragraph and replace the italic "meta description" paragraph's
#    text with the AI-image prompt (keeping its italic run formatting).
# -----------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs($lastIndex - 1)
$boldTitlePara.Range.Delete()

$newLastIndex = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($newLastIndex)
$italicRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End)
$italicRange.Text = "Prompt: Create a cartoon-style image featuring a happy Maya warrior with glasses to fit the theme of the online slot game " + [char]34 + "1 Left Alive." + [char]34 + " The image should be colorful and depict the Maya warrior holding a weapon and surrounded by zombies. The warrior should be wearing traditional Maya clothing and a headpiece, and the glasses should be prominent. Make sure that the image is eye-catching and will draw in players who enjoy action-packed slot games."
